$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Update the "CT 01" scenario row (row 2): test run result flips from Passed -> Failed,
# and the run's output-data timestamp is refreshed to the latest run.
$ws.Range("C2").Value = "Failed"
$ws.Range("H2").Value = "04_05_2020--23_18_51 161"

# Update the "CT 05" scenario row (row 6): refresh its output-data timestamp too.
$ws.Range("H6").Value = "04_05_2020--23_19_24 375"
